# Update 'F' column (想去人数 / interest count) values across all four worksheets
# per the upstream data refresh (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 149
$ws.Range("F6").Value = 781
$ws.Range("F8").Value = 1120
$ws.Range("F9").Value = 301
$ws.Range("F11").Value = 861
$ws.Range("F12").Value = 655
$ws.Range("F13").Value = 177
$ws.Range("F14").Value = 502
$ws.Range("F17").Value = 162
$ws.Range("F18").Value = 2865
$ws.Range("F19").Value = 2601
$ws.Range("F21").Value = 26
$ws.Range("F23").Value = 310
$ws.Range("F24").Value = 221
$ws.Range("F26").Value = 4965
$ws.Range("F31").Value = 273
$ws.Range("F32").Value = 1058
$ws.Range("F34").Value = 46

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1074
$ws.Range("F5").Value = 1074
$ws.Range("F8").Value = 229
$ws.Range("F10").Value = 324
$ws.Range("F11").Value = 11
$ws.Range("F15").Value = 598
$ws.Range("F25").Value = 304
$ws.Range("F26").Value = 270
$ws.Range("F27").Value = 3856
$ws.Range("F32").Value = 38
$ws.Range("F34").Value = 147
$ws.Range("F36").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value = 1271
$ws.Range("F10").Value = 340

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 1271
$ws.Range("F8").Value = 340
$ws.Range("F10").Value = 149
$ws.Range("F12").Value = 781
$ws.Range("F15").Value = 1120
$ws.Range("F16").Value = 301
$ws.Range("F17").Value = 655
$ws.Range("F18").Value = 1074
$ws.Range("F19").Value = 177
$ws.Range("F20").Value = 502
$ws.Range("F22").Value = 162
$ws.Range("F23").Value = 2865
$ws.Range("F24").Value = 2601
$ws.Range("F25").Value = 229
$ws.Range("F26").Value = 310
$ws.Range("F27").Value = 324
$ws.Range("F28").Value = 11
$ws.Range("F29").Value = 221
$ws.Range("F31").Value = 4968
$ws.Range("F34").Value = 598
$ws.Range("F35").Value = 598
$ws.Range("F38").Value = 273
$ws.Range("F42").Value = 304
$ws.Range("F43").Value = 270
$ws.Range("F44").Value = 1058
$ws.Range("F46").Value = 38
$ws.Range("F48").Value = 147
$ws.Range("F50").Value = 46
